$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = New-Object 'object[,]' 24,14
$data[0,0] = 10.33691749721583
$data[0,1] = 4.170258155099462
$data[0,2] = 14.79079604026649
$data[0,3] = 16.10822624963416
$data[0,4] = 0
$data[0,5] = 36.25429821332657
$data[0,6] = 16.24273167720445
$data[0,7] = 0
$data[0,8] = 9.379690435276443
$data[0,9] = 9.729424499044162
$data[0,10] = 0
$data[0,11] = 16.37441099100724
$data[0,12] = 20.03992741058664
$data[0,13] = 25.70930970980926
$data[1,0] = 10.10480197286005
$data[1,1] = 4.015903156279899
$data[1,2] = 14.76347697859329
$data[1,3] = 16.10645183958952
$data[1,4] = 0
$data[1,5] = 36.31820723907194
$data[1,6] = 16.28655910130934
$data[1,7] = 0
$data[1,8] = 9.397094410548208
$data[1,9] = 9.578106636755992
$data[1,10] = 0
$data[1,11] = 16.32018063545015
$data[1,12] = 20.09920009966916
$data[1,13] = 25.77750679315353
$data[2,0] = 9.961493214180592
$data[2,1] = 3.917319327312368
$data[2,2] = 14.74963876133536
$data[2,3] = 16.10828715584081
$data[2,4] = 0
$data[2,5] = 36.36656274402373
$data[2,6] = 16.31574502335022
$data[2,7] = 0
$data[2,8] = 9.408695777678485
$data[2,9] = 9.485564043353678
$data[2,10] = 0
$data[2,11] = 16.28947046972246
$data[2,12] = 20.13729637201298
$data[2,13] = 25.82402941720584
$data[3,0] = 9.902982280304533
$data[3,1] = 3.876225849872186
$data[3,2] = 14.74474241737844
$data[3,3] = 16.10977199788083
$data[3,4] = 0
$data[3,5] = 36.38855490154317
$data[3,6] = 16.32821091247882
$data[3,7] = 0
$data[3,8] = 9.413653915406901
$data[3,9] = 9.447990748812321
$data[3,10] = 0
$data[3,11] = 16.27761566035815
$data[3,12] = 20.15325029836046
$data[3,13] = 25.84415562416693
$data[4,0] = 9.893262448635772
$data[4,1] = 3.869347961754273
$data[4,2] = 14.7439743719748
$data[4,3] = 16.11006310021688
$data[4,4] = 0
$data[4,5] = 36.39234463916965
$data[4,6] = 16.33031543891077
$data[4,7] = 0
$data[4,8] = 9.414491141094098
$data[4,9] = 9.441761455698861
$data[4,10] = 0
$data[4,11] = 16.27568728633338
$data[4,12] = 20.15592540640752
$data[4,13] = 25.84756806299846
$data[5,0] = 9.960704449578046
$data[5,1] = 3.916768796681422
$data[5,2] = 14.74956971413941
$data[5,3] = 16.10830419541983
$data[5,4] = 0
$data[5,5] = 36.36685008581939
$data[5,6] = 16.31591082492721
$data[5,7] = 0
$data[5,8] = 9.408761711179828
$data[5,9] = 9.485056694496494
$data[5,10] = 0
$data[5,11] = 16.28930790832533
$data[5,12] = 20.13750979201583
$data[5,13] = 25.82429611938333
$data[6,0] = 10.25709901736279
$data[6,1] = 4.117848608523504
$data[6,2] = 14.78077004672578
$data[6,3] = 16.1070083647948
$data[6,4] = 0
$data[6,5] = 36.27443932312184
$data[6,6] = 16.25737113931234
$data[6,7] = 0
$data[6,8] = 9.385501568271895
$data[6,9] = 9.677201539010156
$data[6,10] = 0
$data[6,11] = 16.355181030705
$data[6,12] = 20.06001201667788
$data[6,13] = 25.73185836550092
$data[7,0] = 10.82841258691324
$data[7,1] = 4.480492546997845
$data[7,2] = 14.86502004246977
$data[7,3] = 16.12759737335171
$data[7,4] = 0
$data[7,5] = 36.16574600182818
$data[7,6] = 16.16062521468277
$data[7,7] = 0
$data[7,8] = 9.347135672981754
$data[7,9] = 10.05489991613768
$data[7,10] = 0
$data[7,11] = 16.50445794557015
$data[7,12] = 19.92149029132502
$data[7,13] = 25.58753218658116
$data[8,0] = 11.23748544830491
$data[8,1] = 4.725966213498171
$data[8,2] = 14.94063831228132
$data[8,3] = 16.15670328038593
$data[8,4] = 0
$data[8,5] = 36.13032705961124
$data[8,6] = 16.10054087902052
$data[8,7] = 0
$data[8,8] = 9.323345996871147
$data[8,9] = 10.33030528782416
$data[8,10] = 0
$data[8,11] = 16.62577040840469
$data[8,12] = 19.82783564714785
$data[8,13] = 25.50408712141415
$data[9,0] = 11.42037384963327
$data[9,1] = 4.832810881229221
$data[9,2] = 14.97793315498053
$data[9,3] = 16.17294668110715
$data[9,4] = 0
$data[9,5] = 36.12389356584374
$data[9,6] = 16.07559242943749
$data[9,7] = 0
$data[9,8] = 9.313474046056227
$data[9,9] = 10.45460328839771
$data[9,10] = 0
$data[9,11] = 16.68334305829903
$data[9,12] = 19.78697460639721
$data[9,13] = 25.47104421187314
$data[10,0] = 11.48910069973303
$data[10,1] = 4.872557778678487
$data[10,2] = 14.99246403075995
$data[10,3] = 16.17952602533799
$data[10,4] = 0
$data[10,5] = 36.12284999730003
$data[10,6] = 16.06648781015172
$data[10,7] = 0
$data[10,8] = 9.309872069534965
$data[10,9] = 10.50148717868127
$data[10,10] = 0
$data[10,11] = 16.70547466931058
$data[10,12] = 19.77175091890984
$data[10,13] = 25.45923961455691
$data[11,0] = 11.47432371313249
$data[11,1] = 4.864029557103793
$data[11,2] = 14.98931653251744
$data[11,3] = 16.17809005577212
$data[11,4] = 0
$data[11,5] = 36.12301280645324
$data[11,6] = 16.06843340651344
$data[11,7] = 0
$data[11,8] = 9.310641761776758
$data[11,9] = 10.49139879663494
$data[11,10] = 0
$data[11,11] = 16.70069376405381
$data[11,12] = 19.77501853359897
$data[11,13] = 25.46175044075787
$data[12,0] = 11.42603902745602
$data[12,1] = 4.836095270367446
$data[12,2] = 14.97912048868132
$data[12,3] = 16.17347940487608
$data[12,4] = 0
$data[12,5] = 36.12377979849092
$data[12,6] = 16.07483651656982
$data[12,7] = 0
$data[12,8] = 9.313174979040541
$data[12,9] = 10.45846440540825
$data[12,10] = 0
$data[12,11] = 16.68515730058978
$data[12,12] = 19.78571715050598
$data[12,13] = 25.47005884504126
$data[13,0] = 11.39639241362934
$data[13,1] = 4.818891349932453
$data[13,2] = 14.97292801301954
$data[13,3] = 16.17071091796405
$data[13,4] = 0
$data[13,5] = 36.12443097698428
$data[13,6] = 16.0788032502186
$data[13,7] = 0
$data[13,8] = 9.31474439116017
$data[13,9] = 10.43826578163745
$data[13,10] = 0
$data[13,11] = 16.67568335830012
$data[13,12] = 19.79230282145749
$data[13,13] = 25.47524021794295
$data[14,0] = 11.22546286317024
$data[14,1] = 4.718884900094737
$data[14,2] = 14.93825858113851
$data[14,3] = 16.15570187228476
$data[14,4] = 0
$data[14,5] = 36.13094236102113
$data[14,6] = 16.10221929301765
$data[14,7] = 0
$data[14,8] = 9.324010244333529
$data[14,9] = 10.32215854633415
$data[14,10] = 0
$data[14,11] = 16.62205481973384
$data[14,12] = 19.83054098765604
$data[14,13] = 25.50634557961706
$data[15,0] = 11.11973283161407
$data[15,1] = 4.656284318916414
$data[15,2] = 14.91772595830673
$data[15,3] = 16.14726106551762
$data[15,4] = 0
$data[15,5] = 36.13741677071211
$data[15,6] = 16.11719491309541
$data[15,7] = 0
$data[15,8] = 9.329937670230425
$data[15,9] = 10.25064674733823
$data[15,10] = 0
$data[15,11] = 16.58975761405856
$data[15,12] = 19.85444446285873
$data[15,13] = 25.52668775653385
$data[16,0] = 11.05862223057593
$data[16,1] = 4.619825169151828
$data[16,2] = 14.90618922540696
$data[16,3] = 16.14268902387149
$data[16,4] = 0
$data[16,5] = 36.14205170067458
$data[16,6] = 16.1260329097043
$data[16,7] = 0
$data[16,8] = 9.333436414396415
$data[16,9] = 10.20942480206745
$data[16,10] = 0
$data[16,11] = 16.5714066101642
$data[16,12] = 19.8683572288705
$data[16,13] = 25.53885076524803
$data[17,0] = 11.0378823128726
$data[17,1] = 4.607403569745768
$data[17,2] = 14.90233023216792
$data[17,3] = 16.14118969610073
$data[17,4] = 0
$data[17,5] = 36.14377743475179
$data[17,6] = 16.12906384795301
$data[17,7] = 0
$data[17,8] = 9.334636400673018
$data[17,9] = 10.19545360954893
$data[17,10] = 0
$data[17,11] = 16.56523239727233
$data[17,12] = 19.87309607461558
$data[17,13] = 25.5430483915237
$data[18,0] = 11.13101927587559
$data[18,1] = 4.662995285693916
$data[18,2] = 14.91988348559541
$data[18,3] = 16.1481303502862
$data[18,4] = 0
$data[18,5] = 36.13663326214554
$data[18,6] = 16.11557750673867
$data[18,7] = 0
$data[18,8] = 9.329297430546756
$data[18,9] = 10.25826896076962
$data[18,10] = 0
$data[18,11] = 16.59317246263828
$data[18,12] = 19.85188291958797
$data[18,13] = 25.52447439878279
$data[19,0] = 11.44023629082961
$data[19,1] = 4.844319735086024
$data[19,2] = 14.98210430603379
$data[19,3] = 16.17482207001004
$data[19,4] = 0
$data[19,5] = 36.12351671605965
$data[19,6] = 16.07294646426111
$data[19,7] = 0
$data[19,8] = 9.312427214563835
$data[19,9] = 10.46814339397354
$data[19,10] = 0
$data[19,11] = 16.68971188487788
$data[19,12] = 19.78256794491811
$data[19,13] = 25.46759924188366
$data[20,0] = 11.63921555722518
$data[20,1] = 4.958663554571794
$data[20,2] = 15.02514387935627
$data[20,3] = 16.19476135528181
$data[20,4] = 0
$data[20,5] = 36.12306146166339
$data[20,6] = 16.04708279341243
$data[20,7] = 0
$data[20,8] = 9.302195969636482
$data[20,9] = 10.6042109313705
$data[20,10] = 0
$data[20,11] = 16.75472362052973
$data[20,12] = 19.73872037328383
$data[20,13] = 25.4345553970269
$data[21,0] = 11.53332275379228
$data[21,1] = 4.898022649155172
$data[21,2] = 15.00195843787193
$data[21,3] = 16.18389236827607
$data[21,4] = 0
$data[21,5] = 36.12256168830083
$data[21,6] = 16.06070390013038
$data[21,7] = 0
$data[21,8] = 9.307583989901893
$data[21,9] = 10.53170341621075
$data[21,10] = 0
$data[21,11] = 16.71985463259264
$data[21,12] = 19.76198999758644
$data[21,13] = 25.45181355953889
$data[22,0] = 11.12591768143526
$data[22,1] = 4.659962717719563
$data[22,2] = 14.91890723289738
$data[22,3] = 16.14773647206589
$data[22,4] = 0
$data[22,5] = 36.13698464341999
$data[22,6] = 16.11630802496043
$data[22,7] = 0
$data[22,8] = 9.329586599464825
$data[22,9] = 10.2548232938224
$data[22,10] = 0
$data[22,11] = 16.59162793238892
$data[22,12] = 19.85304046260532
$data[22,13] = 25.52547360003674
$data[23,0] = 10.67541275587771
$data[23,1] = 4.38596381572502
$data[23,2] = 14.83979172402793
$data[23,3] = 16.11956239977939
$data[23,4] = 0
$data[23,5] = 36.18735947123493
$data[23,6] = 16.18486605728315
$data[23,7] = 0
$data[23,8] = 9.356740897625166
$data[23,9] = 9.952902402158394
$data[23,10] = 0
$data[23,11] = 16.46198385856963
$data[23,12] = 19.95753265838665
$data[23,13] = 25.62261356220631

$ws.Range("B2:O25").Value = $data
